# The workbook was edited to hold the records in a custom bean:
# a simple 2-row x 4-column table of string data was added to Sheet1,
# with the first cell (A1) formatted with a date number format (mmm-yy).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header-ish) values
$ws.Range("A1").Value = "one"
$ws.Range("B1").Value = "onetwo"
$ws.Range("C1").Value = "one3"
$ws.Range("D1").Value = "one4"

# Row 2 values
$ws.Range("A2").Value = "two1"
$ws.Range("B2").Value = "two2"
$ws.Range("C2").Value = "two3"
$ws.Range("D2").Value = "two4"

# A1 carries a custom number format (numFmtId 17 -> "mmm-yy")
$ws.Range("A1").NumberFormat = "mmm-yy"

# Match the saved selection/active cell (D2) from the source file
$null = $ws.Range("D2").Select()
